# Apply the authored changes:
#  1. Refresh the cached "datetimeFigureOut" footer-date field text from
#     1/19/2022 to 1/20/2022 across the slide master and every slide layout.
#  2. On slide 2, move/resize the "Subscription Options" textbox and
#     shorten its caption to "Subscriptions".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text (slide master + all custom layouts)
# ---------------------------------------------------------------------
$oldDate = "1/19/2022"
$newDate = "1/20/2022"

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) "Subscription Options" -> "Subscriptions" textbox on slide 2
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$grp = $s2.Shapes.Item(1)

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "Subscription Options") {
        $sh.Left = 536.0424499511719
        $sh.Top = 115.60496139526367
        $sh.Width = 133.13213348388672
        $sh.Height = 31.504723548889164
        $sh.TextFrame.TextRange.Text = "Subscriptions"
    }
}
